# Signal R done For Admin Dashboard
# Update a few existing order rows to "Accepted" status, then append new
# incoming order rows (simulating new real-time submissions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: mark them as Accepted ---
$ws.Range("D7").Value = "Accepted"
$ws.Range("D9").Value = "Accepted"
$ws.Range("D10").Value = "Accepted"

# --- Append new rows with new submissions ---
$newRows = @(
    @("Farid Abdull", "eqarayev4@std.beu.edu.az", 66, "Pending"),
    @("Farid Abdull", "eqarayev4@std.beu.edu.az", 156, "Rejected"),
    @("Farid Abdull", "eqarayev4@std.beu.edu.az", 66, "Accepted"),
    @("Farid Abdull", "eqarayev4@std.beu.edu.az", 156, "Accepted"),
    @("Elmar Qarayev", "elmarqarayev69@gmail.com", 27, "Pending"),
    @("Elmar Qarayev", "elmarqarayev69@gmail.com", 9, "Accepted"),
    @("Elmar Qarayev", "elmarqarayev69@gmail.com", 18, "Accepted")
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
}
